$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data through row 343 (9 Aug 2021).
# This update ("aggiornamento al 23 agosto 2021") appends 14 more daily
# rows (344-357, 10-23 Aug 2021) with the same layout:
#   A = date (serial), B = nuovi positivi, C = somma mobile 7gg.,
#   D = somma mobile 7gg. per 100mila abitanti.
#
# First, stamp column A of every new row with the date-cell formatting
# (style index 2: bold, centered, bordered, custom date format) by
# copying the format from the last existing date cell, A343.
for ($r = 344; $r -le 357; $r++) {
    $ws.Range("A343").Copy($ws.Range("A$r"))
}

# Now fill in the actual values for the new rows.
$ws.Range("A344").Value = 44418
$ws.Range("B344").Value = 1
$ws.Range("C344").Value = 5
$ws.Range("D344").Value = 60.76810889645115

$ws.Range("A345").Value = 44419
$ws.Range("B345").Value = 0
$ws.Range("C345").Value = 5
$ws.Range("D345").Value = 60.76810889645115

$ws.Range("A346").Value = 44420
$ws.Range("B346").Value = 0
$ws.Range("C346").Value = 5
$ws.Range("D346").Value = 60.76810889645115

$ws.Range("A347").Value = 44421
$ws.Range("B347").Value = 2
$ws.Range("C347").Value = 6
$ws.Range("D347").Value = 72.92173067574137

$ws.Range("A348").Value = 44422
$ws.Range("B348").Value = 0
$ws.Range("C348").Value = 4
$ws.Range("D348").Value = 48.61448711716091

$ws.Range("A349").Value = 44423
$ws.Range("B349").Value = 0
$ws.Range("C349").Value = 3
$ws.Range("D349").Value = 36.46086533787069

$ws.Range("A350").Value = 44424
$ws.Range("B350").Value = 0
$ws.Range("C350").Value = 3
$ws.Range("D350").Value = 36.46086533787069

$ws.Range("A351").Value = 44425
$ws.Range("B351").Value = 0
$ws.Range("C351").Value = 2
$ws.Range("D351").Value = 24.30724355858046

$ws.Range("A352").Value = 44426
$ws.Range("B352").Value = 0
$ws.Range("C352").Value = 2
$ws.Range("D352").Value = 24.30724355858046

$ws.Range("A353").Value = 44427
$ws.Range("B353").Value = 0
$ws.Range("C353").Value = 2
$ws.Range("D353").Value = 24.30724355858046

$ws.Range("A354").Value = 44428
$ws.Range("B354").Value = 1
$ws.Range("C354").Value = 1
$ws.Range("D354").Value = 12.15362177929023

$ws.Range("A355").Value = 44429
$ws.Range("B355").Value = 0
$ws.Range("C355").Value = 1
$ws.Range("D355").Value = 12.15362177929023

$ws.Range("A356").Value = 44430
$ws.Range("B356").Value = 0
$ws.Range("C356").Value = 1
$ws.Range("D356").Value = 12.15362177929023

$ws.Range("A357").Value = 44431
$ws.Range("B357").Value = 0
$ws.Range("C357").Value = 1
$ws.Range("D357").Value = 12.15362177929023
